# Uppdaterat estimering av projektet
# Rework the "Backlog items" table on Blad2 (sheet2) to match the new estimation/plan.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the backlog rows (row 3 header stays the same) ---

# Row 4: Test Identity (unchanged)
$ws.Cells.Item(4,2).Value = 1
$ws.Cells.Item(4,3).Value = "Test Identity"
$ws.Cells.Item(4,4).Value = 2
$ws.Cells.Item(4,5).Value = 1
$ws.Cells.Item(4,6).Value = "ToDo"
$ws.Cells.Item(4,7).Value = "Köra med en test implementering av Identity med Google"

# Row 5: Implementera Identity
$ws.Cells.Item(5,2).Value = 2
$ws.Cells.Item(5,3).Value = "Implementera Identity"
$ws.Cells.Item(5,4).Value = 4
$ws.Cells.Item(5,5).Value = 1
$ws.Cells.Item(5,6).Value = "ToDo"
$ws.Cells.Item(5,7).Value = "Implementera Identity, få igång alla funktioner, plus layout"

# Row 6: Layout för login-sida (new item)
$ws.Cells.Item(6,2).Value = 3
$ws.Cells.Item(6,3).Value = "Layout för login-sida"
$ws.Cells.Item(6,4).Value = 2
$ws.Cells.Item(6,5).Value = 1
$ws.Cells.Item(6,6).Value = "ToDo"
$ws.Cells.Item(6,7).Value = "Implementera all layout för login sida och dess respektive partial views"

# Row 7: Modulera Databas
$ws.Cells.Item(7,2).Value = 3
$ws.Cells.Item(7,3).Value = "Modulera Databas"
$ws.Cells.Item(7,4).Value = 2
$ws.Cells.Item(7,5).Value = 1
$ws.Cells.Item(7,6).Value = "ToDo"
$ws.Cells.Item(7,7).Value = "Finslipa, modulera slutgiltlig databas"

# Row 8: Implementera rest databas
$ws.Cells.Item(8,2).Value = 4
$ws.Cells.Item(8,3).Value = "Implementera rest databas"
$ws.Cells.Item(8,4).Value = 2
$ws.Cells.Item(8,5).Value = 1
$ws.Cells.Item(8,6).Value = "ToDo"
$ws.Cells.Item(8,7).Value = "Implementering av resterande databas med entity framework code first"

# Row 9: Implementera ToDo
$ws.Cells.Item(9,2).Value = 6
$ws.Cells.Item(9,3).Value = "Implementera ToDo"
$ws.Cells.Item(9,4).Value = 12
$ws.Cells.Item(9,5).Value = 1
$ws.Cells.Item(9,6).Value = "ToDo"
$ws.Cells.Item(9,7).Value = "Lägga till alla dess respektive funktioner"

# Row 10: only an item number + placeholder description remain ("…")
$ws.Cells.Item(10,2).Value = 7
$ws.Cells.Item(10,3).Value = "…"
$ws.Range("D10:G10").ClearContents()

# --- Remove the old rows 11-15 entirely ---
$ws.Rows("11:15").Delete()

# --- Column widths: col G narrower, new col H added ---
$ws.Columns("G").ColumnWidth = 65.30729166666667
$ws.Columns("H").ColumnWidth = 37.166666666666664

# --- Selection moves to C11 ---
[void]$ws.Range("C11").Select()
